# Adding Parameterization and automating the third test
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddCustomerTest")

# Make sure this sheet is the active one (this also flips tabSelected
# off the previously-active "OpenAccountTest" sheet and sets the
# workbook's active-tab back to the first sheet).
$ws.Activate()

# Row 3: new parameterized customer
$ws.Range("A3").Value = "Rahul"
$ws.Range("B3").Value = "Arora"
$ws.Range("C3").Value = "A234wd"
$ws.Range("D3").Value = "Customer added successfully"

# Row 4: new parameterized customer
$ws.Range("A4").Value = "Ishita"
$ws.Range("B4").Value = "Arora"
$ws.Range("C4").Value = "A234wd"
$ws.Range("D4").Value = "Customer added successfully"

# Row 5: new parameterized customer (third, automated test)
$ws.Range("A5").Value = "Rohit"
$ws.Range("B5").Value = "Sehgal"
$ws.Range("C5").Value = "A234wd"
$ws.Range("D5").Value = "Customer added successfully"

# Widen column D so the long alert text fits (best-fit autosize)
$ws.Columns.Item(4).ColumnWidth = 26.3

# Leave the cursor on the last entered cell
$ws.Range("B5").Select()
